$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the "_old" / "_new" header suffixes to the new version tags
#    ("_FV2410" / "_FV2504") on the header row (row 1, columns A:U).
# ---------------------------------------------------------------------------
$oldHeaders = @(
  "Segmentname_old",
  "Segmentgruppe_old",
  "Segment_old",
  "Datenelement_old",
  "Segment ID_old",
  "Code_old",
  "Qualifier_old",
  "Beschreibung_old",
  "Bedingungsausdruck_old",
  "Bedingung_old"
)
$newHeaders = @(
  "Segmentname_new",
  "Segmentgruppe_new",
  "Segment_new",
  "Datenelement_new",
  "Segment ID_new",
  "Code_new",
  "Qualifier_new",
  "Beschreibung_new",
  "Bedingungsausdruck_new",
  "Bedingung_new"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
  $col = $i + 1
  $ws.Cells.Item(1, $col).Value = ($oldHeaders[$i] -replace "_old$", "_FV2410")
}

# column 11 ("diff") is left untouched

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
  $col = $i + 12
  $ws.Cells.Item(1, $col).Value = ($newHeaders[$i] -replace "_new$", "_FV2504")
}

# ---------------------------------------------------------------------------
# 2) Turn the used range into an Excel Table ("Table1") so the regenerated
#    AHB diff sheet gets header-row filtering again.
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U73")
$lo = $ws.ListObjects.Add(
  [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
  $tableRange,
  $null,
  [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split/freeze at row 2) and keep row 1 visible
#    while scrolling.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
